$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.911620497703552
$ws.Range("B1").Value = 1.89979076385498
$ws.Range("C1").Value = 1.807050466537476
$ws.Range("D1").Value = 1.003948211669922
$ws.Range("E1").Value = 0.6872348785400391
